$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the strategy labels / acronyms (column B, rows 2-8)
$ws.Range("B2").Value = "DB Search"
$ws.Range("B3").Value = "SB Search (BS*FS)"
$ws.Range("B4").Value = "DB Search + BS*FS"
$ws.Range("B5").Value = "Scopus + BS*FS"
$ws.Range("B6").Value = "Scopus + BS||FS"
$ws.Range("B7").Value = "Scopus + BS+FS"
$ws.Range("B8").Value = "Scopus + FS+BS"

# 2. Swap the numeric data (columns C:N) between rows 3 and 4
#    (the order of the two strategies was swapped). Literal values are
#    written directly (rather than read back and reassigned) to avoid any
#    loss of floating point precision during the swap.
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 489
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.02249488752556237
$ws.Range("J3").Value = 0.7857142857142857
$ws.Range("K3").Value = 0.0437375745526839
$ws.Range("L3").Value = 0.02249488752556237
$ws.Range("M3").Value = 0.7857142857142857
$ws.Range("N3").Value = 0.0437375745526839

$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 932
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.01502145922746781
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0.02959830866807611
$ws.Range("L4").Value = 0.01502145922746781
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 0.02959830866807611
